# ---------------------------------------------------------------------------
# "ire keys in one sheet"
#
# - Insert a new worksheet "input_ire_mus" between "input_1" and "all_keys"
#   containing the NET / NET_IR_UI rows (entity/key/min_cap/max_cap) that
#   used to live scattered inside "input".
# - Move "all_keys" to the end of the tab strip (after "config").
# - Refresh the autofilter / filter-database on "input" (now A1:D1504) and
#   tidy up its selection.
# - Drop the stray tabSelected on "input_1".
# - Bump the counter on "config" (B4: 4 -> 6) and make it the active tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsInput   = $wb.Worksheets.Item("input")
$wsInput1  = $wb.Worksheets.Item("input_1")
$wsAllKeys = $wb.Worksheets.Item("all_keys")
$wsConfig  = $wb.Worksheets.Item("config")

# ---------------------------------------------------------------------------
# 1) Create "input_ire_mus" right after "input_1"
# ---------------------------------------------------------------------------
$wsIre = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsInput1)
$wsIre.Name = "input_ire_mus"

$rows = @(
  @("220KV-KORBA-BUDIPADAR","NET",-10000,10000),
  @("220KV-RAIGARH-BUDIPADAR","NET",-10000,10000),
  @("400KV-Raigarh-JHARSUGUDA","NET",-10000,10000),
  @("400KV-SIPAT-RANCHI","NET",-10000,10000),
  @("765KV-DHARJAYGARH-JHARSUGUDA","NET",-10000,10000),
  @("765KV-DHARJAYGARH-RANCHI","NET",-10000,10000),
  @("WR-ER","NET_IR_UI",-10000,10000),
  @("220KV-Badod-Kota","NET",-10000,10000),
  @("220KV-Badod-Modak","NET",-10000,10000),
  @("220KV-Malanpur-Auriya","NET",-10000,10000),
  @("220KV-Mehgaon-Auriya","NET",-10000,10000),
  @("400KV-KANSARI-BHINMAL","NET",-10000,10000),
  @("400KV-KANSARI-KANKROLI","NET",-10000,10000),
  @("400KV-Sujalpur-RAPP","NET",-10000,10000),
  @("400KV-VSTPS-RIHAND","NET",-10000,10000),
  @("765KV-GWALIOR-AGRA","NET",-10000,10000),
  @("765KV-GWALIOR-Jaipur","NET",-10000,10000),
  @("765KV-GWALIOR-Orai","NET",-10000,10000),
  @("765KV-JABALPUR-Orai","NET",-10000,10000),
  @("765KV-SATNA-Orai","NET",-10000,10000),
  @("HVDC400KV-Vindyachal(PS)-RIHAND","NET",-10000,10000),
  @("HVDC500KV-Mundra-Mohindargarh","NET",-10000,10000),
  @("HVDC800KV-CHAMPA-KURUKSHETRA","NET",-10000,10000),
  @("WR-SR","NET_IR_UI",-10000,10000),
  @("220KV-XELDEM-AMBEWADI","NET",-10000,10000),
  @("400KV-Kolhapur GIS-Narendra Kudgi","NET",-10000,10000),
  @("765KV-Solapur-Raichur","NET",-10000,10000),
  @("765KV-Wardha-Nizamabad","NET",-10000,10000),
  @("HVDC500KV-BHADRAWATI-RAMAGUNDAM","NET",-10000,10000),
  @("WR-NR","NET_IR_UI",-10000,10000)
)

$data = New-Object 'object[,]' $rows.Length,4
for ($i = 0; $i -lt $rows.Length; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

$lastRow = $rows.Length + 1
$wsIre.Range("A2:D$lastRow").Value = $data

# Header row, text + formatting copied from the "input" header so the
# bold/centered/bordered style (and shared strings) line up exactly.
$hdr = New-Object 'object[,]' 1,4
$hdr[0,0] = "entity"
$hdr[0,1] = "key"
$hdr[0,2] = "min_cap"
$hdr[0,3] = "max_cap"
$wsIre.Range("A1:D1").Value = $hdr

$wsInput.Range("A1:D1").Copy()
$wsIre.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Column widths (bestFit-ish) matching the source sheets as closely as the
# engine's character-width rounding allows.
$wsIre.Columns.Item(1).ColumnWidth = 38.83333333333333
$wsIre.Columns.Item(2).ColumnWidth = 11.0
$wsIre.Columns.Item(3).ColumnWidth = 7.666666666666667
$wsIre.Columns.Item(4).ColumnWidth = 7.833333333333333

$wsIre.Range("C15").Select()

# ---------------------------------------------------------------------------
# 2) Push "all_keys" to the end of the tab strip (after "config")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAllKeys.Move([System.Reflection.Missing]::Value, $lastSheet)

# ---------------------------------------------------------------------------
# 3) "input": drop the stale topLeftCell/selection, re-point the selection,
#    and (re)apply the autofilter over the full A1:D1504 range.
# ---------------------------------------------------------------------------
$wsInput.Range("B9").Select()
$wsInput.Range("A1:D1504").AutoFilter()
$wb.Names.Item("input!_FilterDatabase").RefersTo = "=input!`$A`$1:`$D`$1504"

# ---------------------------------------------------------------------------
# 4) "input_1": selection stays at B12 - just touch it so any stray
#    tabSelected flag is dropped once "config" becomes the active sheet.
# ---------------------------------------------------------------------------
$wsInput1.Range("B12").Select()

# ---------------------------------------------------------------------------
# 5) "config": bump the counter, keep selection at B4, make it the active tab.
# ---------------------------------------------------------------------------
$wsConfig.Range("B4").Value = 6
$wsConfig.Range("B4").Select()
$wsConfig.Activate()
